$d = $word.ActiveDocument

$replacements = @(
    @{old = "824÷9="; new = "266÷7="},
    @{old = "300÷7="; new = "847÷8="},
    @{old = "685÷3="; new = "625÷5="},
    @{old = "389÷8="; new = "434÷8="},
    @{old = "612÷3="; new = "582÷7="},
    @{old = "487÷8="; new = "485÷5="},
    @{old = "940÷3="; new = "847÷4="},
    @{old = "294÷6="; new = "317÷8="},
    @{old = "965÷5="; new = "349÷5="},
    @{old = "169÷3="; new = "247÷4="},
    @{old = "777÷8="; new = "436÷8="},
    @{old = "895÷5="; new = "267÷8="},
    @{old = "606÷7="; new = "461÷7="},
    @{old = "774÷3="; new = "264÷2="},
    @{old = "676÷7="; new = "469÷6="},
    @{old = "385÷4="; new = "649÷2="},
    @{old = "339÷3="; new = "193÷6="},
    @{old = "985÷6="; new = "691÷6="},
    @{old = "570÷4="; new = "483÷6="},
    @{old = "140÷4="; new = "426÷5="},
    @{old = "896÷6="; new = "860÷9="},
    @{old = "327÷8="; new = "219÷2="},
    @{old = "377÷5="; new = "262÷6="},
    @{old = "329÷8="; new = "658÷9="},
    @{old = "352÷6="; new = "877÷5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
